# edit.ps1 — apply the "Append: 2025-11-06 01:20 JST" scrape update to 案件情報.xlsx
#
# A new job listing was scraped and inserted into the sorted (by priority score)
# list on row 4, pushing the previously-4th..12th rows down by one. Every row's
# "retrieved at" timestamp (column A) is refreshed to the time of this run, and
# the hyperlinks on column F are rebuilt to point at the correct (shifted) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert a new row at position 4 to make room for the new listing,
# shifting the existing rows 4-12 down to 5-13.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new job listing.
$ws.Range("A4").Value = "2025-11-06 01:20:17"
$ws.Range("B4").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("G4").Value = 310
$ws.Range("H4").Value = "🔥AI,Ai"

# Refresh the "retrieved at" timestamp for every data row (2-13) to this run's time.
$newTimestamp = "2025-11-06 01:20:17"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Rebuild the hyperlinks on column F: the row insert shifted the cell contents
# but left the old hyperlink map (ref -> rId) stale, so start clean and re-add
# one hyperlink per data row, in order, restoring the Hyperlink cell style.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5416301")
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5420440")
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5427956")
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5416328")
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5427011")
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5427648")
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5427397")
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5427338")
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5427682")
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5427793")
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5427459")
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5427699")
$ws.Range("F13").Style = "Hyperlink"

# Dimension / used-range bookkeeping is handled automatically by the engine,
# but make sure the sheet thinks its extent is A1:H13 explicitly as well.
$ws.Range("A1:H13").Select()

